$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date formatting (style) from A4 into A5, then set the new value
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)
$ws.Cells.Item(5, 1).Value = 42587.819618055553

$ws.Cells.Item(5, 2).Value = "Random"
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 79
$ws.Cells.Item(5, 9).Value = 21
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 91
$ws.Cells.Item(5, 13).Value = 9
